$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$orangeOld = "📙"
$blueOld   = "📘"
$orangeNew = "+3"
$blueNew   = "⚠️"

# Scan every used cell for the two emoji markers and replace them in place,
# keeping the literal string on the left of -eq so boolean-valued cells
# (TRUE/FALSE) are never coerced into matching a non-empty string.
$used = $ws.UsedRange
foreach ($cell in $used.Cells) {
    $v = $cell.Value2
    if ($orangeOld -eq $v) {
        # Force text formatting so "+3" is stored as a literal string
        # rather than being auto-converted to the number 3.
        $cell.NumberFormat = "@"
        $cell.Value = $orangeNew
    } elseif ($blueOld -eq $v) {
        $cell.Value = $blueNew
    }
}
